$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138. This shifts the existing rows
# 138-181 down to become rows 139-182 (preserving all their values and
# formatting), and leaves a blank (but correctly formatted) row 138
# ready to be filled with the new data record.
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new data record.
$ws.Range("A138").Value = 4
$ws.Range("B138").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C138").Value = "Los Lagos"
$ws.Range("D138").Value = 45173
$ws.Range("E138").Value = 10
$ws.Range("F138").Value = "Fruta"
$ws.Range("G138").Value = 100104
$ws.Range("H138").Value = "Frutos de pepita"
$ws.Range("I138").Value = 100104003
$ws.Range("J138").Value = "Membrillo"
$ws.Range("K138").Value = "Champion"
$ws.Range("L138").Value = "Primera"
$ws.Range("M138").Value = 50
$ws.Range("N138").Value = 15000
$ws.Range("O138").Value = 15000
$ws.Range("P138").Value = 15000
$ws.Range("Q138").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R138").Value = "Región de O'Higgins"
$ws.Range("S138").Value = 833
$ws.Range("T138").Value = 18
